$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '64.429.12'
Set-TextValue "E2" '  +1.54%  '
Set-TextValue "D3" '3.360.88'
Set-TextValue "E3" '  +1.39%  '
Set-TextValue "E4" '  +0.05%  '
Set-TextValue "D5" '555.72'
Set-TextValue "E5" '  +1.34%  '
Set-TextValue "D6" '173.69'
Set-TextValue "E6" '  +0.52%  '
Set-TextValue "E7" '  +1.56%  '
Set-TextValue "D8" '3.352.32'
Set-TextValue "E8" '  +1.38%  '
Set-TextValue "E9" '  +0.07%  '
Set-TextValue "E10" '  +7.55%  '
Set-TextValue "E11" '  +2.78%  '
Set-TextValue "D12" '53.64'
Set-TextValue "E12" '  +1.13%  '
Set-TextValue "D13" '0.0000273'
Set-TextValue "E13" '  +3.19%  '
Set-TextValue "E14" '  +1.68%  '
Set-TextValue "D15" '3.899.85'
Set-TextValue "E15" '  +1.53%  '
Set-TextValue "E16" '  -0.08%  '
Set-TextValue "E17" '  +1.73%  '
Set-TextValue "D18" '3.363.02'
Set-TextValue "E18" '  +1.71%  '
Set-TextValue "D19" '64.456.83'
Set-TextValue "E19" '  +1.73%  '
Set-TextValue "D20" '11.69'
Set-TextValue "E20" '  +0.18%  '
Set-TextValue "D21" '0.985'
Set-TextValue "E21" '  +1.81%  '
Set-TextValue "D22" '458.62'
Set-TextValue "E22" '  +8.22%  '
Set-TextValue "D23" '4.87'
Set-TextValue "E23" '  +9.41%  '
Set-TextValue "D24" '4.09'
Set-TextValue "E24" '  +1.00%  '
Set-TextValue "D25" '86.00'
Set-TextValue "E25" '  +3.54%  '
Set-TextValue "D26" '13.61'
Set-TextValue "E26" '  +2.41%  '
Set-TextValue "D27" '2.96'
Set-TextValue "E27" '  +8.77%  '
Set-TextValue "E28" '  +1.15%  '
Set-TextValue "E29" '  +0.33%  '
Set-TextValue "E30" '  +4.49%  '
Set-TextValue "E31" '  +3.70%  '
Set-TextValue "E32" '  +0.49%  '
Set-TextValue "D33" '571.13'
Set-TextValue "E33" '  -0.80%  '
Set-TextValue "D34" '61.13'
Set-TextValue "E34" '  +5.13%  '
Set-TextValue "E35" '  +1.20%  '
Set-TextValue "D37" '3.62'
Set-TextValue "E37" '  +3.95%  '
Set-TextValue "E38" '  -4.87%  '
Set-TextValue "D39" '35.18'
Set-TextValue "E39" '  +0.26%  '
Set-TextValue "E40" '  +0.10%  '
Set-TextValue "E41" '  +1.00%  '
Set-TextValue "E42" '  +0.13%  '
Set-TextValue "D43" '3.069.65'
Set-TextValue "E43" '  -1.68%  '
Set-TextValue "D44" '2.80'
Set-TextValue "E44" '  +0.89%  '
Set-TextValue "E45" '  +2.81%  '
Set-TextValue "E47" '  +0.68%  '
Set-TextValue "D48" '3.16'
Set-TextValue "E48" '  -0.57%  '
Set-TextValue "E49" '  -0.89%  '
Set-TextValue "D50" '137.84'
Set-TextValue "E50" '  +3.42%  '
Set-TextValue "E51" '  +1.69%  '
